$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update raw timing measurements for the "CAE" run (columns F and G) ---
# Column F: CAE SMC2 dataset timings
$ws.Range("F2").Value = 144.86000000000001
$ws.Range("F3").Value = 78.7
$ws.Range("F4").Value = 56.73
$ws.Range("F5").Value = 46.45
$ws.Range("F6").Value = 39.36
$ws.Range("F7").Value = 35.020000000000003

# Column G: CAE HEA dataset timings
$ws.Range("G2").Value = 709.24
$ws.Range("G3").Value = 680.24
$ws.Range("G4").Value = 669.65
$ws.Range("G5").Value = 658.48
$ws.Range("G6").Value = 656
$ws.Range("G7").Value = 650

# AVERAGE/STDEV in J6:K7 and the normalized Z-scores in F9:G14 are formulas
# that recompute automatically from the new F2:G7 inputs.

# --- Restore the window/selection UI state recorded at save time ---
$win = $excel.ActiveWindow
$win.Left = 28680

$ws.Range("G79").Select()
